$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 76
$ws.Range("H76").Value = 2923.7083
$ws.Range("I76").Value = 2868
$ws.Range("K76").Value = 2868
$ws.Range("M76").Value = -2553

# Row 79
$ws.Range("H79").Value = 2923.7083
$ws.Range("I79").Value = 2868
$ws.Range("K79").Value = 2868
$ws.Range("M79").Value = -1776

# Row 105
$ws.Range("H105").Value = 39600
$ws.Range("J105").Value = 39600
$ws.Range("L105").Value = 39600
$ws.Range("N105").Value = -46588

# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents() | Out-Null

# Row 121
$ws.Range("H121").Value = 716.6667
$ws.Range("J121").Value = 1000
$ws.Range("L121").Value = 3000
$ws.Range("N121").Value = -6494

# Row 135
$ws.Range("H135").Value = 1614.05
$ws.Range("I135").Value = 1220.2941
$ws.Range("K135").Value = 10982.6469
$ws.Range("M135").Value = -8447.6469

# Row 138
$ws.Range("H138").Value = 3386.012
$ws.Range("I138").Value = 2179.96
$ws.Range("J138").Value = 3905.862
$ws.Range("K138").Value = 6539.88
$ws.Range("L138").Value = 11717.586
$ws.Range("M138").Value = -1399.88
$ws.Range("N138").Value = -21997.586

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 88
$ws.Range("H88").Value = 2380
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 2633.3333
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 2633.3333
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -3445.3333

# Row 91
$ws.Range("H91").Value = 2380
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 2633.3333
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 2633.3333
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -5441.3333

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents() | Out-Null

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Range("H86").Value = 21090
$ws.Range("I86").Value = 1366.2354
$ws.Range("J86").Value = 58346
$ws.Range("K86").Value = 1366.2354
$ws.Range("L86").Value = 58346
$ws.Range("M86").Value = -243.2354
$ws.Range("N86").Value = -60592

# Row 89
$ws.Range("H89").Value = 21090
$ws.Range("I89").Value = 1366.2354
$ws.Range("J89").Value = 58346
$ws.Range("K89").Value = 6831.177
$ws.Range("L89").Value = 291730
$ws.Range("M89").Value = -1215.177
$ws.Range("N89").Value = -302962

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents() | Out-Null

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 62
$ws.Range("H62").Value = 4032.4
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4032.4
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4032.4
$ws.Range("M62").ClearContents() | Out-Null
$ws.Range("N62").Value = -5280.4

# Row 65
$ws.Range("H65").Value = 4032.4
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4032.4
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 20162
$ws.Range("M65").ClearContents() | Out-Null
$ws.Range("N65").Value = -26402

# Row 132
$ws.Range("H132").Value = 2645.5454
$ws.Range("I132").Value = 2031.9546
$ws.Range("J132").Value = 3872.7273
$ws.Range("K132").Value = 6095.8638
$ws.Range("L132").Value = 11618.1819
$ws.Range("M132").Value = -3565.8638
$ws.Range("N132").Value = -16678.1819

# Row 134
$ws.Range("H134").Value = 1578.0571
$ws.Range("I134").Value = 756.7406999999999
$ws.Range("J134").Value = 4350
$ws.Range("K134").Value = 2270.2221
$ws.Range("L134").Value = 13050
$ws.Range("M134").Value = 264.7779
$ws.Range("N134").Value = -18120

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 17
$ws.Range("H17").Value = 462.13333
$ws.Range("I17").Value = 231.44444
$ws.Range("J17").Value = 808.1667
$ws.Range("K17").Value = 694.33332
$ws.Range("L17").Value = 2424.5001
$ws.Range("M17").Value = -525.33332
$ws.Range("N17").Value = -2762.5001

# Row 34
$ws.Range("H34").Value = 6893.5293
$ws.Range("I34").Value = 223.75
$ws.Range("J34").Value = 12822.223
$ws.Range("K34").Value = 671.25
$ws.Range("L34").Value = 38466.669
$ws.Range("M34").Value = -587.25
$ws.Range("N34").Value = -38634.669

# Row 39
$ws.Range("H39").Value = 1789.2858
$ws.Range("J39").Value = 2355
$ws.Range("L39").Value = 7065
$ws.Range("N39").Value = -7653

# Row 55
$ws.Range("H55").Value = 3085
$ws.Range("I55").Value = 420
$ws.Range("J55").Value = 3917.8125
$ws.Range("K55").Value = 1260
$ws.Range("L55").Value = 11753.4375
$ws.Range("M55").Value = -1083
$ws.Range("N55").Value = -12107.4375

# Row 92
$ws.Range("H92").Value = 4800
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 4800
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 14400
$ws.Range("M92").ClearContents() | Out-Null
$ws.Range("N92").Value = -16896

# Row 120
$ws.Range("H120").Value = 18955.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 25
$ws.Range("H25").Value = 70009
$ws.Range("J25").Value = 70009
$ws.Range("L25").Value = 70009
$ws.Range("N25").Value = -71067

# Row 116
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents() | Out-Null

# Row 126
$ws.Range("H126").Value = 3612.6428
$ws.Range("I126").Value = 1794.25
$ws.Range("K126").Value = 5382.75
$ws.Range("M126").Value = -2912.75

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents() | Out-Null

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 14
$ws.Range("H14").Value = 371432
$ws.Range("J14").Value = 18004
$ws.Range("L14").Value = 18004
$ws.Range("N14").Value = -18348

# Row 132
$ws.Range("H132").Value = 4018.4285
$ws.Range("I132").Value = 2414.8572
$ws.Range("J132").Value = 4820.2144
$ws.Range("K132").Value = 7244.571599999999
$ws.Range("L132").Value = 14460.6432
$ws.Range("M132").Value = -4714.571599999999
$ws.Range("N132").Value = -19520.6432

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents() | Out-Null

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 123
$ws.Range("H123").Value = 29700
$ws.Range("J123").Value = 29700
$ws.Range("L123").Value = 29700
$ws.Range("N123").Value = -39500

# Row 132
$ws.Range("H132").Value = 5122.41
$ws.Range("I132").Value = 2016.8636
$ws.Range("J132").Value = 9141.352999999999
$ws.Range("K132").Value = 6050.5908
$ws.Range("L132").Value = 27424.059
$ws.Range("M132").Value = -3520.5908
$ws.Range("N132").Value = -32484.059

# Row 136
$ws.Range("H136").Value = 1113.4884
$ws.Range("I136").Value = 767.96875
$ws.Range("J136").Value = 2118.6365
$ws.Range("K136").Value = 2303.90625
$ws.Range("L136").Value = 6355.9095
$ws.Range("M136").Value = 246.09375
$ws.Range("N136").Value = -11455.9095

# Row 140
$ws.Range("H140").Value = 34360.75
$ws.Range("J140").Value = 38192.9
$ws.Range("L140").Value = 38192.9
$ws.Range("N140").Value = -48552.9

# Row 141
$ws.Range("H141").Value = 28011.316
$ws.Range("J141").Value = 28011.316
$ws.Range("L141").Value = 28011.316
$ws.Range("N141").Value = -38371.316
